$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.44810000000002
$ws.Range("D5").Value = -8.632499999999995
$ws.Range("D9").Value = -8.503500000000001
$ws.Range("D11").Value = -8.289300000000004
$ws.Range("A21").Value = -21.1534
$ws.Range("D21").Value = -8.166900000000005
$ws.Range("A23").Value = -21.34160000000002
$ws.Range("A25").Value = -22.33530000000003
